$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Characteristics")

# --- Falling Creek Reservoir (row 2): new DIC (mg/L) measurement ---
$ws.Range("I2").Value = 8.3

# --- Lake Sunapee (row 4): new DOC (mg/L) and DIC (mg/L) measurements ---
$ws.Range("H4").Value = 2.36
$ws.Range("I4").Value = 2

# These cells were highlighted yellow as a "missing data" placeholder;
# now that real values have been entered, clear that highlight (xlNone).
$ws.Range("I2").Interior.Pattern = -4142
$ws.Range("H4").Interior.Pattern = -4142
$ws.Range("I4").Interior.Pattern = -4142

# Restore the active selection left by the editor after entering the data.
$ws.Range("J7").Select()

$wb.Save()
